$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidated "Absent" (column H) results for rows whose attendance
# totals (column D) indicate an absence that was previously left blank
# or miscalculated.
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
